$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
# Old layout: StaffId | FromDuration | ToDuration | Reason | StartDate | EndDate | StartTime | EndTime | ApplicationTypeName
# New layout: StaffId | ApplicationType | StartDate | EndDate | StartTime | EndTime | Reason | TotalDays | TotalHours | StartDuration | EndDuration
$ws.Range("A1").Value = "StaffId"
$ws.Range("B1").Value = "ApplicationType"
$ws.Range("C1").Value = "StartDate"
$ws.Range("D1").Value = "EndDate"
$ws.Range("E1").Value = "StartTime"
$ws.Range("F1").Value = "EndTime"
$ws.Range("G1").Value = "Reason"
$ws.Range("H1").Value = "TotalDays"
$ws.Range("I1").Value = "TotalHours"
$ws.Range("J1").Value = "StartDuration"
$ws.Range("K1").Value = "EndDuration"

# --- Sample/format row (row 2) -----------------------------------------
# Date columns get a short-date display format, time/duration columns get
# a duration display format.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").NumberFormat = "mm:ss.0"
$ws.Range("F2").NumberFormat = "mm:ss.0"
$ws.Range("M2").NumberFormat = "mm:ss.0"

# --- Column widths (best-fit, in characters) ------------------------------
$ws.Columns.Item(1).ColumnWidth = 5.833333333333333
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(3).ColumnWidth = 8.5
$ws.Columns.Item(4).ColumnWidth = 7.666666666666667
$ws.Columns.Item(5).ColumnWidth = 8.666666666666666
$ws.Columns.Item(6).ColumnWidth = 7.666666666666667
$ws.Columns.Item(7).ColumnWidth = 6.666666666666667
$ws.Columns.Item(8).ColumnWidth = 8.666666666666666
$ws.Columns.Item(9).ColumnWidth = 9.666666666666666
$ws.Columns.Item(10).ColumnWidth = 11.833333333333334
$ws.Columns.Item(11).ColumnWidth = 11

# --- Selection ------------------------------------------------------------
[void]$ws.Range("C8").Select()
